$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read the current data & per-row "is text" flag out of column A ---
$lastRow = $ws.UsedRange.Rows.Count
$firstDataRow = 2
$n = $lastRow - $firstDataRow + 1

# Template cells that already carry the two body styles used in this sheet:
#   A3  -> text style  (centered, non-bold, Text "@" format)
#   A42 -> number style (centered, non-bold, General format)
# Stash copies of those two styles in a scratch area (far outside the data
# range) before we start clearing/rewriting A2:B<lastRow>, since the originals
# live inside the region we are about to wipe out.
$scratchText = $ws.Cells.Item(1, 26)   # Z1
$scratchNum = $ws.Cells.Item(2, 26)    # Z2
$ws.Range("A3").Copy() | Out-Null
$scratchText.PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$scratchNum.PasteSpecial(-4122) | Out-Null
$styleTextTemplate = $scratchText
$styleNumTemplate = $scratchNum

$vals = @()
$isText = @()
for ($i = 0; $i -lt $n; $i++) {
    $cell = $ws.Cells.Item($firstDataRow + $i, 1)
    $vals += ,($cell.Value2)
    $isText += ,($cell.NumberFormat -eq "@")
}

# The first three rows ("1-400", "1501-1687", "14792-15192") stay put in column A,
# and what used to be their paired single-filter hits move alongside them into
# column B. Everything else (the remaining sorted region list, plus the trailing
# numeric position list) moves entirely into column B, starting at row 2 - the
# single filters are no longer split out into their own column A rows.

$keepInA = 3          # "1-400", "1501-1687", "14792-15192"

# Clear out all of the old body (rows firstDataRow..lastRow, columns A and B)
$ws.Range($ws.Cells.Item($firstDataRow,1), $ws.Cells.Item($lastRow,2)).Clear() | Out-Null

# Re-write column A: the three region labels (same text style as before)
for ($i = 0; $i -lt $keepInA; $i++) {
    $r = $firstDataRow + $i
    $c = $ws.Cells.Item($r, 1)
    $styleTextTemplate.Copy()
    $c.PasteSpecial(-4122) | Out-Null
    $c.Value = $vals[$i]
}

# Re-write column B with every remaining value, starting at row 2
$bCount = $n - $keepInA
for ($i = 0; $i -lt $bCount; $i++) {
    $srcIndex = $keepInA + $i
    $r = $firstDataRow + $i
    $c = $ws.Cells.Item($r, 2)
    if ($isText[$srcIndex]) {
        $styleTextTemplate.Copy()
    } else {
        $styleNumTemplate.Copy()
    }
    $c.PasteSpecial(-4122) | Out-Null
    $c.Value = $vals[$srcIndex]
}

$excel.CutCopyMode = 0
$scratchText.Clear() | Out-Null
$scratchNum.Clear() | Out-Null

$newLastRow = $firstDataRow + $bCount - 1

# Remove the now-unused trailing rows so the sheet dimension / used range shrinks
if ($lastRow -gt $newLastRow) {
    $rowSpec = "$($newLastRow + 1):$lastRow"
    $ws.Rows($rowSpec).Delete() | Out-Null
}

# Update the active selection / view to match the edited sheet
$ws.Range("A5").Select() | Out-Null
